$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing (pre-edit) values we need to re-use, to avoid retyping long literals
$origWebData = $ws.Cells.Item(2,4).Value2      # "TC09_..._WebData.xlsx"
$origNeo4jData = $ws.Cells.Item(2,3).Value2    # "TC09_..._Neo4jData.xlsx"
$origQuery = $ws.Cells.Item(1,1).Value2        # "query"
$origStatQueryLabel = $ws.Cells.Item(1,2).Value2  # "StatQuery"
$origDbExcel = $ws.Cells.Item(1,3).Value2      # "dbExcel"
$origWebExcel = $ws.Cells.Item(1,4).Value2     # "WebExcel"

# --- New long query texts ---
$samplesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN [''Border Collie'']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'
$filesQuery = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Border Collie'']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '''') AS `File Name`,
        coalesce(f.file_type, '''') AS `File Type`,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`,
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'
$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Border Collie'']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`,
        coalesce(co.cohort_description, '''') AS `Cohort`
'
$statQueryText = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Border Collie'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# --- Insert a new column A (shifts old A..D -> B..E) ---
$ws.Columns.Item(1).Insert()

# --- Insert two new rows after row 2 (so we end up with 4 data rows total) ---
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# --- Column A labels (written first, in row order, so shared-string indices
#     for the label text match the order Excel originally produced them) ---
$ws.Cells.Item(1,1).Value2 = "TabName"
$ws.Cells.Item(2,1).Value2 = "CasesTab"
$ws.Cells.Item(3,1).Value2 = "SamplesTab"
$ws.Cells.Item(4,1).Value2 = "FilesTab"

# --- Query column (B): Samples, then Files, then Cases (matches source order) ---
$ws.Cells.Item(3,2).Value2 = $samplesQuery
$ws.Cells.Item(4,2).Value2 = $filesQuery
$ws.Cells.Item(2,2).Value2 = $casesQuery

# --- Remaining columns (C, D, E) for each row ---
$ws.Cells.Item(1,2).Value2 = $origQuery
$ws.Cells.Item(1,3).Value2 = $origStatQueryLabel
$ws.Cells.Item(1,4).Value2 = $origDbExcel
$ws.Cells.Item(1,5).Value2 = $origWebExcel

$ws.Cells.Item(2,3).Value2 = $statQueryText
$ws.Cells.Item(2,4).Value2 = $origNeo4jData
$ws.Cells.Item(2,5).Value2 = $origWebData

$ws.Cells.Item(3,3).Value2 = $statQueryText
$ws.Cells.Item(3,4).Value2 = $origNeo4jData
$ws.Cells.Item(3,5).Value2 = $origWebData

$ws.Cells.Item(4,3).Value2 = $statQueryText
$ws.Cells.Item(4,4).Value2 = $origNeo4jData
$ws.Cells.Item(4,5).Value2 = $origWebData

# --- Wrap-text style for the long-text cells (style index 1 = wrapText) ---
$ws.Cells.Item(2,2).WrapText = $true
$ws.Cells.Item(2,3).WrapText = $true
$ws.Cells.Item(3,2).WrapText = $true
$ws.Cells.Item(3,3).WrapText = $true
$ws.Cells.Item(4,2).WrapText = $true
$ws.Cells.Item(4,3).WrapText = $true

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 10.90625
$ws.Columns.Item(2).ColumnWidth = 75.81640625
$ws.Columns.Item(3).ColumnWidth = 75.81640625
$ws.Columns.Item(4).ColumnWidth = 70.26953125
$ws.Columns.Item(5).ColumnWidth = 28.54296875

# --- Row heights (headless runtime has no text-measuring AutoFit, so set explicitly) ---
$ws.Rows.Item(2).RowHeight = 275.5
$ws.Rows.Item(3).RowHeight = 232
$ws.Rows.Item(4).RowHeight = 246.5

# --- Sheet view: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 40
$ws.Range("C2").Select()
